# Insert a new weekly record at row 25 (Vega Monumental Concepción - Poroto verde).
# This pushes the existing rows 25-72 down to 26-73, preserving all of their data,
# and fills the newly created row 25 with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 25..72 down by one row.
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the new weekly entry.
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44979
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112031
$ws.Range("G25").Value = "Poroto verde"
$ws.Range("H25").Value = "Magnum"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 26000
$ws.Range("M25").Value = 25500
$ws.Range("N25").Value = "$/saco 25 kilos"
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 1020
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"
